$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (data for MuSCs sending cluster -> removed)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Itgb6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.835941000000001
$ws.Range("H2").Value = 17.507823
$ws.Range("I2").Value = 0.03643643319117328
$ws.Range("J2").Value = 0.03643643319117327
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.6957970000000001
$ws.Range("N2").Value = 2.087391
$ws.Range("O2").Value = 0.9232770860517062
$ws.Range("P2").Value = 0.9232770860517063
$ws.Range("Q2").Value = 4.060630239977002
$ws.Range("R2").Value = 36.54567215979301
$ws.Range("S2").Value = 0.03364092386286413
$ws.Range("T2").Value = 0.03364092386286413

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Itgb6"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.835941000000001
$ws.Range("H3").Value = 17.507823
$ws.Range("I3").Value = 0.03643643319117328
$ws.Range("J3").Value = 0.03643643319117327
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.05781966666666667
$ws.Range("N3").Value = 0.173459
$ws.Range("O3").Value = 0.07672291394829377
$ws.Range("P3").Value = 0.07672291394829378
$ws.Range("Q3").Value = 0.3374321633063334
$ws.Range("R3").Value = 3.036889469757
$ws.Range("S3").Value = 0.002795509328309142
$ws.Range("T3").Value = 0.002795509328309142

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Itgb6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.50798033333334
$ws.Range("H4").Value = 52.52394100000001
$ws.Range("I4").Value = 0.1093102818770573
$ws.Range("J4").Value = 0.1093102818770573
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.6957970000000001
$ws.Range("N4").Value = 2.087391
$ws.Range("O4").Value = 0.9232770860517062
$ws.Range("P4").Value = 0.9232770860517063
$ws.Range("Q4").Value = 12.18200019199234
$ws.Range("R4").Value = 109.638001727931
$ws.Range("S4").Value = 0.1009236785269401
$ws.Range("T4").Value = 0.1009236785269401

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Itgb6"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.50798033333334
$ws.Range("H5").Value = 52.52394100000001
$ws.Range("I5").Value = 0.1093102818770573
$ws.Range("J5").Value = 0.1093102818770573
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05781966666666667
$ws.Range("N5").Value = 0.173459
$ws.Range("O5").Value = 0.07672291394829377
$ws.Range("P5").Value = 0.07672291394829378
$ws.Range("Q5").Value = 1.012305586879889
$ws.Range("R5").Value = 9.110750281919001
$ws.Range("S5").Value = 0.008386603350117202
$ws.Range("T5").Value = 0.008386603350117204

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Itgb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 136.8238143333333
$ws.Range("H6").Value = 410.471443
$ws.Range("I6").Value = 0.8542532849317694
$ws.Range("J6").Value = 0.8542532849317694
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6957970000000001
$ws.Range("N6").Value = 2.087391
$ws.Range("O6").Value = 0.9232770860517062
$ws.Range("P6").Value = 0.9232770860517063
$ws.Range("Q6").Value = 95.20159954169034
$ws.Range("R6").Value = 856.8143958752131
$ws.Range("S6").Value = 0.7887124836619021
$ws.Range("T6").Value = 0.7887124836619021

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Itgb6"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 136.8238143333333
$ws.Range("H7").Value = 410.471443
$ws.Range("I7").Value = 0.8542532849317694
$ws.Range("J7").Value = 0.8542532849317694
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.05781966666666667
$ws.Range("N7").Value = 0.173459
$ws.Range("O7").Value = 0.07672291394829377
$ws.Range("P7").Value = 0.07672291394829378
$ws.Range("Q7").Value = 7.911107336815222
$ws.Range("R7").Value = 71.199966031337
$ws.Range("S7").Value = 0.06554080126986743
$ws.Range("T7").Value = 0.06554080126986743
